$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 16 (pushing the
# existing row 16 and everything below it down by one row, growing the
# table from A1:R137 to A1:R138).
$ws.Rows.Item(16).Insert()

# The new record is a duplicate of the (now shifted-down) old row 16 in
# every column except the date, so seed it by copying that row...
$ws.Range("A17:R17").Copy()
$ws.Range("A16:R16").PasteSpecial()

# ...then correct the date (Fecha) to the new reporting date, 2021-12-31.
$ws.Range("D16").Value = 44561
